$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enter the new quiz scores for "Li, Meng Yong" (currently row 34) ---
$ws.Range("C34").Value = 9
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 1
$ws.Range("I34").Value = 1
$ws.Range("J34").Value = 1
$ws.Range("K34").Value = 1
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0
$ws.Range("R34").Value = 0
$ws.Range("S34").Value = 0
$ws.Range("T34").Value = 0
$ws.Range("U34").Value = 0

# --- Mark the roster "counted" flag (column A) for the two rows whose flag
#     changed alongside this update ---
$ws.Range("A28").Value = 1   # Wu, Nina
$ws.Range("A29").Value = 1   # Cai, Yu row slot -> flag now also set
$ws.Range("A38").ClearContents()   # Yang, Shuo loses its flag

# --- Re-sort the roster (name + scores, not the flag column) by total
#     score descending, then by name ascending, same as the workbook's
#     existing sortState ---
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add2($ws.Range("V2:V39"), 0, 2, $null, 0) | Out-Null
$sort.SortFields.Add2($ws.Range("B2:B39"), 0, 1, $null, 0) | Out-Null
$sort.SetRange($ws.Range("B2:V39"))
$sort.Header = 2
$sort.Apply()
